{"js": "// Change \"Represent this game in normal form\" to\n// \"Represent this game in extensive form\" in the homework question\n// about the driver/accident/insurance game.\n\nconst body = context.document.body;\n\n// Find the run of text containing \"normal form\" inside the target\n// sentence and replace just that phrase, preserving the rest of the\n// sentence (and the run's formatting).\nconst searchResults = body.search(\"normal form\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find text \"normal form\" in the document.');\n}\n\n// Replace from the last match to the first so earlier matches' anchors\n// aren't disturbed by edits made to later ones.\nfor (let i = searchResults.items.length - 1; i >= 0; i--) {\n  searchResults.items[i].insertText(\"extensive form\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Change \"Represent this game in normal form\" to\n# \"Represent this game in extensive form\" in the homework question\n# about the driver/accident/insurance game.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"normal form\"\n$find.Replacement.Text = \"extensive form\"\n$find.Forward = $true\n$find.Wrap = 1                 # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2 -- replace every occurrence of \"normal form\" in the\n# document body with \"extensive form\".\n$find.Execute(\"normal form\", $false, $true, $false, $false, $false, $true, 1, $false, \"extensive form\", 2) | Out-Null\n"}
